# Mark the "table 3" rows for years 1999-2006 as "Finished" (column H) and
# add notes (column I) describing open questions, matching the upstream
# commit "Finished up to 2006_tb3.xlsx. Continue with 2007_tb3".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("I1").Value = "Notes"

# Mark rows as Finished ("X") in column H
$ws.Range("H9").Value  = "X"
$ws.Range("H11").Value = "X"
$ws.Range("H13").Value = "X"
$ws.Range("H15").Value = "X"
$ws.Range("H17").Value = "X"
$ws.Range("H19").Value = "X"
$ws.Range("H21").Value = "X"
$ws.Range("H23").Value = "X"

# Add notes in column I (order matters for shared-string table indices:
# the short "Where do other expenses come from?" note is used first,
# before the longer combined note)
$ws.Range("I15").Value = "Where do other expenses come from?"
$ws.Range("I17").Value = "Where do other expenses come from?"
$ws.Range("I19").Value = "Where do other expenses come from?"
$ws.Range("I21").Value = "Where do other expenses come from?"
$ws.Range("I23").Value = "Where do other expenses come from?"
$ws.Range("I13").Value = "Where do other expenses come from and why are only some amphibians counted in subtotal?"

# Move/restore the active selection to match the final state of the file
$ws.Range("I23").Select() | Out-Null
